# Planning V7.0 - "Master - Changes of implications"
#
# Two textual edits to Task 10 ("Address any relevant implications...") content:
#
# 1. Expand on why no copyrighted / illegal / explicit images are a concern:
#    "No copyrighted images. No illegal or explicit images etc." ->
#    "No copyrighted images. No illegal or explicit images etc, this should come
#     from the pages earlier so it should be fine."
#
# 2. Add a closing sentence about efficiency / user control to the paragraph that
#    explains why the page was made similar to the showcase page:
#    "...which make the user feel "comfortable" returning a dog." ->
#    "...which make the user feel "comfortable" returning a dog. This allows for
#     the program to be efficient to use and give the user freedom and control
#     over the program."

$d = $word.ActiveDocument

$curlyOpen  = [char]8220   # U+201C LEFT DOUBLE QUOTATION MARK
$curlyClose = [char]8221   # U+201D RIGHT DOUBLE QUOTATION MARK

# --- Edit 1 -----------------------------------------------------------------
$old1 = "No copyrighted images. No illegal or explicit images etc."
$new1 = "No copyrighted images. No illegal or explicit images etc, this should come from the pages earlier so it should be fine."

$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) {
    Write-Output "WARNING: edit 1 target text not found"
}

# --- Edit 2 -----------------------------------------------------------------
$old2 = "the user feel " + $curlyOpen + "comfortable" + $curlyClose + " returning a dog."
$new2 = "the user feel " + $curlyOpen + "comfortable" + $curlyClose + " returning a dog. This allows for the program to be efficient to use and give the user freedom and control over the program."

$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $found2) {
    Write-Output "WARNING: edit 2 target text not found"
}

Write-Output "Edit 1 applied: $found1"
Write-Output "Edit 2 applied: $found2"
